# Generate Report for Handoff
#
# - Overview!E2/F2, zh-cn!C2, de-de!C2: status text changes from
#   "Handed back: in sync with en-US" to "Ready for handoff".
# - Overview!G2: Latest HO Xliff Generate Date refreshed.
# - zh-cn!H2: Latest Handoff Datetime refreshed.
# - de-de!H2: Latest Handoff Datetime refreshed (same stamp as Overview!G2).
# - The Status / Xliff-date columns (Overview E & F, zh-cn C, de-de C) are
#   narrowed now that the status text is shorter.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff"
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# --- Refreshed timestamps
$wsOverview.Range("G2").Value = "2016-08-16 00:54:28"
$wsZhCn.Range("H2").Value = "2016-08-16 00:54:23"
$wsDeDe.Range("H2").Value = "2016-08-16 00:54:28"

# --- Narrow the columns that used to hold the longer status text.
# The sheet's raw column width (OOXML "width") is ColumnWidth*6+5 pixels,
# rounded to the nearest pixel; 16.3333... is the ColumnWidth that lands
# closest to the target 17.216-character raw width.
$newWidth = 16.333333333333332
$wsOverview.Range("E1").ColumnWidth = $newWidth
$wsOverview.Range("F1").ColumnWidth = $newWidth
$wsZhCn.Range("C1").ColumnWidth = $newWidth
$wsDeDe.Range("C1").ColumnWidth = $newWidth
